$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (col index 2), rows 2-25
$BValues = @(0.6111739799903546,0.5687642829294077,0.5429188198673671,0.5324360502573597,0.5306984002389754,0.5427772444618597,0.5965111530578042,0.7034031555612899,0.7828429776635915,0.819175228334359,0.8329607982139748,0.8299906238260064,0.8203088313983926,0.8143819942825132,0.7804724406332753,0.7597194195236625,0.7478011938877671,0.7437690614929124,0.7619267183089562,0.8231518719019277,0.8633250742964833,0.8418695465833821,0.7609287581027502,0.674325425920415)
for ($i = 0; $i -lt $BValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $BValues[$i]
}

# Column C (col index 3), rows 2-25
$CValues = @(0.1864339052119419,0.1867681177701606,0.1869959736744136,0.1870945466779546,0.1871112608505037,0.18699727987158,0.186544455889571,0.1858351020294116,0.1854214326345556,0.1852562934352733,0.1851970491288029,0.1852096624744561,0.1852513535594582,0.1852773183719592,0.1854326865454894,0.185533884829141,0.1855942614929447,0.1856150774133312,0.1855228876701318,0.1852390187596455,0.1850726621147416,0.1851597034081749,0.1855278526409379,0.1860080237745052)
for ($i = 0; $i -lt $CValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 3).Value = $CValues[$i]
}

# Column D (col index 4), rows 2-25
$DValues = @(0.05180973644298348,0.05010185833353376,0.04904136516626068,0.04860624778928724,0.04853381887907915,0.0490355089686787,0.05122333189697059,0.05541885827135928,0.0584427635924456,0.05980555961131984,0.06031975711225357,0.06020909875352487,0.05984790043843446,0.05962641261127288,0.05835344216558269,0.05756922209833704,0.0571169566854266,0.05696362154498757,0.05765282832063434,0.0599540439127253,0.06144714676264584,0.06065125260880677,0.05761503432784565,0.05429408431543692)
for ($i = 0; $i -lt $DValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 4).Value = $DValues[$i]
}

# Column E (col index 5), rows 2-25
$EValues = @(0.1194009475103428,0.1189125529881245,0.1186737650074363,0.1185918329572502,0.1185791573128263,0.1186725977665688,0.1192198786790577,0.1207773541942601,0.1222165488624398,0.1229352700255539,0.1232166276965394,0.1231556236105256,0.1229582333127652,0.1228385229406044,0.122170865963021,0.1217776689311378,0.1215575394536934,0.1214840429436812,0.1218189017433211,0.123015962221686,0.1238518914585853,0.1234008419602048,0.1218002419333963,0.1203042032760493)
for ($i = 0; $i -lt $EValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 5).Value = $EValues[$i]
}

# Column F (col index 6), rows 2-25
$FValues = @(1.068408293977036,1.069010116001749,1.069978404952749,1.070523567438656,1.070623186081505,1.069985147510678,1.06849153236017,1.070314351611025,1.074553548128563,1.077112569984678,1.078172310049027,1.077940042181766,1.077197937707027,1.076755188781178,1.074399000504968,1.073115076265665,1.07243595023975,1.072216203310745,1.073245609234874,1.077413450248827,1.080666016076222,1.078881684840354,1.073186411438598,1.069312189369448)
for ($i = 0; $i -lt $FValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 6).Value = $FValues[$i]
}

# Column K (col index 11), rows 2-25
$KValues = @(0.3152865101614566,0.2791099968643493,0.2568984416937923,0.2478477703461408,0.2463449706154677,0.2567763776385732,0.3028129194118208,0.3930818299906491,0.4593818799429812,0.4895362296731207,0.5009536692641916,0.4984947879212314,0.4904755794527205,0.4855633906392427,0.4574110601098766,0.4401386463702579,0.4302034883363035,0.4268395404741909,0.4419773821011574,0.492831056309825,0.5260587187786712,0.5083254050420294,0.4411461052686434,0.368664207379851)
for ($i = 0; $i -lt $KValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 11).Value = $KValues[$i]
}

# Column L (col index 12), rows 2-25
$LValues = @(0.1920005324637515,0.1849179742812481,0.1806710816900647,0.1789661052878415,0.1786845472901746,0.1806479837753869,0.1895373665811206,0.2077757740710808,0.2216666099973423,0.2280925395396736,0.2305412088001333,0.230013164228339,0.2282936865903764,0.2272424486488944,0.2212488062662743,0.2175992386418528,0.2155101714930936,0.2148045803687495,0.2179866998598357,0.2287983238894213,0.2359535483608255,0.2321265300237911,0.2178115002383265,0.2027555266339078)
for ($i = 0; $i -lt $LValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 12).Value = $LValues[$i]
}

# Column O (col index 15), rows 2-25
$OValues = @(3.848243232169608,3.865075140854231,3.877423588115761,3.882961979803781,3.883912203790487,3.877496230868786,3.853628928863486,3.822806458634261,3.80991382820838,3.806168782541278,3.805055604291169,3.805281780235987,3.80607108738397,3.8065942839375,3.810201246253939,3.812957061874471,3.814741657846554,3.8153801550732,3.812643049516197,3.805830970203999,3.803156663553494,3.804421273950823,3.812784390768627,3.829432809506415)
for ($i = 0; $i -lt $OValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 15).Value = $OValues[$i]
}

